# Update "想去人数" (interest count) values in column F for sheets
# "展览" and "全部类型", reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Updates common to both the "展览" and "全部类型" sheets.
$commonUpdates = @{
    3  = 105
    4  = 288
    6  = 595
    8  = 2067
    13 = 285
    16 = 131
    19 = 84
    20 = 3370
    21 = 81
    22 = 523
    24 = 18
    25 = 84
    29 = 64
    30 = 206
    32 = 648
    33 = 2035
    34 = 383
}

# Sheet "展览" (sheet1): row 11 -> 4506
$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $commonUpdates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
$ws1.Cells.Item(11, 6).Value = 4506

# Sheet "全部类型" (sheet4): row 11 -> 4507
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $commonUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $commonUpdates[$row]
}
$ws4.Cells.Item(11, 6).Value = 4507
